# Fix std (X = fg_mf_total) computation issue and correct row-data shuffling
# that resulted from re-deriving compound property rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ row=2; A="palmitic acid, tms derivative"; B="hexadecanoic acid"; C="C16H32O2"; D="CCCCCCCCCCCCCCCC(=O)O"; E=256.42; F=6.4; G="palmitic acid"; H=16; I=32; J=2; K=0.7494579205990172; L=0.125793619842446; M=0.1247874580765931; N=15; O=0; P=0; Q=1; R=0; S=0.8244793697839481; T=0; U=0; V=0.1755596287341081; W=0; X=0.9998867256637168 },
    @{ row=3; A="9,12-octadecadienoic acid (z,z)-, tms derivative"; B="(9z,12z)-octadeca-9,12-dienoic acid"; C="C18H32O2"; D="CCCCCC=CCC=CCCCCCCCC(=O)O"; E=280.4; F=6.8; G="9,12-octadecadienoic acid (z,z)-"; H=18; I=32; J=2; K=0.7710342368045648; L=0.1150356633380885; M=0.1141155492154066; N=17; O=0; P=0; Q=1; R=0; S=0.8396398002853066; T=0; U=0; V=0.1605456490727532; W=0; X=0.9998867256637168 },
    @{ row=4; A="benzoic acid, deriv."; B="benzoic acid"; C="C7H6O2"; D="C1=CC=C(C=C1)C(=O)O"; E=122.12; F=1.9; G="benzoic acid"; H=7; I=6; J=2; K=0.6884785456927611; L=0.04952505732066819; M=0.2620209629872257; N=0; O=6; P=0; Q=1; R=0; S=0; T=0.6313953488372093; U=0; V=0.3686292171634458; W=0; X=0.9998867256637168 },
    @{ row=5; A="myristic acid, tms derivative"; B="tetradecanoic acid"; C="C14H28O2"; D="CCCCCCCCCCCCCC(=O)O"; E=228.37; F=5.3; G="myristic acid"; H=14; I=28; J=2; K=0.7363226343214958; L=0.1235889127293427; M=0.1401147261023777; N=13; O=0; P=0; Q=1; R=0; S=0.8029031834303979; T=0; U=0; V=0.1971230897228183; W=0; X=0.9998867256637168 },
    @{ row=6; A="benzene-1,2-diol, deriv."; B="benzene-1,2-diol"; C="C6H6O2"; D="C1=CC=C(C(=C1)O)O"; E=110.11; F=0.9; G="benzene-1,2-diol"; H=6; I=6; J=2; K=0.6544909635818728; L=0.05492689129052766; M=0.2906003087821269; N=0; O=6; P=2; Q=0; R=0; S=0; T=0.6911088911088911; U=0.3089092725456362; V=0; W=0; X=0.9998867256637168 },
    @{ row=7; A="4-oxopentanoic acid, deriv."; B="4-oxopentanoic acid"; C="C5H8O3"; D="CC(=O)CCC(=O)O"; E=116.11; F=-0.5; G="4-oxopentanoic acid"; H=5; I=8; J=3; K=0.5172250452157436; L=0.06945138230987856; M=0.413375247610025; N=1; O=0; P=0; Q=1; R=1; S=0.1208078546206184; T=0; U=0; V=0.3877099302385669; W=0.491533890276462; X=0.9998867256637168 },
    @{ row=8; A="phenol, deriv."; B="phenol"; C="C6H6O"; D="C1=CC=C(C=C1)O"; E=94.11; F=1.5; G="phenol"; H=6; I=6; J=1; K=0.765763468281798; L=0.06426522154925088; M=0.1700031877590054; N=0; O=6; P=1; Q=0; R=0; S=0; T=0.8193178195728402; U=0.1807140580172139; V=0; W=0; X=0.9998867256637168 },
    @{ row=9; A="palmitelaidic acid, tms derivative"; B="(e)-hexadec-9-enoic acid"; C="C16H30O2"; D="CCCCCCC=CCCCCCCCC(=O)O"; E=254.41; F=6.4; G="palmitelaidic acid"; H=16; I=30; J=2; K=0.7553791124562713; L=0.1188632522306513; M=0.1257733579654888; N=15; O=0; P=0; Q=1; R=0; S=0.8230690617507173; T=0; U=0; V=0.1769466609016941; W=0; X=0.9998867256637168 },
    @{ row=10; A="9-octadecenoic acid, (z)-, tms derivative"; B="(z)-octadec-9-enoic acid"; C="C18H34O2"; D="CCCCCCCCC=CCCCCCCCC(=O)O"; E=282.5; F=6.5; G="9-octadecenoic acid, (z)-"; H=18; I=34; J=2; K=0.7653026548672566; L=0.121316814159292; M=0.1132672566371681; N=17; O=0; P=0; Q=1; R=0; S=0.8405345132743363; T=0; U=0; V=0.1593522123893805; W=0; X=0.9998867256637168 },
    @{ row=11; A="hexadecanoic acid, deriv."; B="hexadecanoic acid"; C="C16H32O2"; D="CCCCCCCCCCCCCCCC(=O)O"; E=256.42; F=6.4; G="hexadecanoic acid"; H=16; I=32; J=2; K=0.7494579205990172; L=0.125793619842446; M=0.1247874580765931; N=15; O=0; P=0; Q=1; R=0; S=0.8244793697839481; T=0; U=0; V=0.1755596287341081; W=0; X=0.9998867256637168 },
    @{ row=12; A="(9z,12z)-octadeca-9,12-dienoic acid, deriv."; B="(9z,12z)-octadeca-9,12-dienoic acid"; C="C18H32O2"; D="CCCCCC=CCC=CCCCCCCCC(=O)O"; E=280.4; F=6.8; G="(9z,12z)-octadeca-9,12-dienoic acid"; H=18; I=32; J=2; K=0.7710342368045648; L=0.1150356633380885; M=0.1141155492154066; N=17; O=0; P=0; Q=1; R=0; S=0.8396398002853066; T=0; U=0; V=0.1605456490727532; W=0; X=0.9998867256637168 },
    @{ row=13; A="9-octadecenoic acid, (e)-, deriv."; B="(e)-octadec-9-enoic acid"; C="C18H34O2"; D="CCCCCCCCC=CCCCCCCCC(=O)O"; E=282.5; F=6.5; G="9-octadecenoic acid, (e)-"; H=18; I=34; J=2; K=0.7653026548672566; L=0.121316814159292; M=0.1132672566371681; N=17; O=0; P=0; Q=1; R=0; S=0.8405345132743363; T=0; U=0; V=0.1593522123893805; W=0; X=0.9998867256637168 }
)

foreach ($r in $rows) {
    $ws.Range("A$($r.row)").Value = $r.A
    $ws.Range("B$($r.row)").Value = $r.B
    $ws.Range("C$($r.row)").Value = $r.C
    $ws.Range("D$($r.row)").Value = $r.D
    $ws.Range("E$($r.row)").Value = $r.E
    $ws.Range("F$($r.row)").Value = $r.F
    $ws.Range("G$($r.row)").Value = $r.G
    $ws.Range("H$($r.row)").Value = $r.H
    $ws.Range("I$($r.row)").Value = $r.I
    $ws.Range("J$($r.row)").Value = $r.J
    $ws.Range("K$($r.row)").Value = $r.K
    $ws.Range("L$($r.row)").Value = $r.L
    $ws.Range("M$($r.row)").Value = $r.M
    $ws.Range("N$($r.row)").Value = $r.N
    $ws.Range("O$($r.row)").Value = $r.O
    $ws.Range("P$($r.row)").Value = $r.P
    $ws.Range("Q$($r.row)").Value = $r.Q
    $ws.Range("R$($r.row)").Value = $r.R
    $ws.Range("S$($r.row)").Value = $r.S
    $ws.Range("T$($r.row)").Value = $r.T
    $ws.Range("U$($r.row)").Value = $r.U
    $ws.Range("V$($r.row)").Value = $r.V
    $ws.Range("W$($r.row)").Value = $r.W
    $ws.Range("X$($r.row)").Value = $r.X
}
